$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-14 Monday", "2025-04-15 Tuesday"),
    @("869÷6=144, 5", "261÷5=52, 1"),
    @("236÷6=39, 2", "770÷5=154, 0"),
    @("355÷8=44, 3", "613÷6=102, 1"),
    @("994÷2=497, 0", "137÷7=19, 4"),
    @("533÷5=106, 3", "451÷7=64, 3"),
    @("328÷7=46, 6", "805÷2=402, 1"),
    @("195÷5=39, 0", "187÷8=23, 3"),
    @("948÷9=105, 3", "393÷7=56, 1"),
    @("920÷4=230, 0", "764÷3=254, 2"),
    @("458÷2=229, 0", "831÷7=118, 5"),
    @("452÷2=226, 0", "507÷4=126, 3"),
    @("900÷2=450, 0", "805÷9=89, 4"),
    @("429÷8=53, 5", "913÷3=304, 1"),
    @("107÷3=35, 2", "313÷3=104, 1"),
    @("216÷7=30, 6", "589÷2=294, 1"),
    @("190÷5=38, 0", "789÷3=263, 0"),
    @("410÷6=68, 2", "115÷7=16, 3"),
    @("721÷7=103, 0", "100÷5=20, 0"),
    @("762÷9=84, 6", "893÷8=111, 5"),
    @("238÷5=47, 3", "349÷2=174, 1"),
    @("679÷8=84, 7", "697÷6=116, 1"),
    @("287÷5=57, 2", "826÷3=275, 1"),
    @("752÷9=83, 5", "456÷6=76, 0"),
    @("239÷8=29, 7", "341÷2=170, 1"),
    @("443÷8=55, 3", "414÷7=59, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
